$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-01-28 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-29 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("26+28=54", $true, $false, $false, $false, $false, $true, 1, $false, "34-32=2", 2) | Out-Null
$d.Content.Find.Execute("62+31=93", $true, $false, $false, $false, $false, $true, 1, $false, "66+27=93", 2) | Out-Null
$d.Content.Find.Execute("21+61=82", $true, $false, $false, $false, $false, $true, 1, $false, "11+41=52", 2) | Out-Null
$d.Content.Find.Execute("11+27=38", $true, $false, $false, $false, $false, $true, 1, $false, "91-64=27", 2) | Out-Null
$d.Content.Find.Execute("42-12=30", $true, $false, $false, $false, $false, $true, 1, $false, "81+0=81", 2) | Out-Null
$d.Content.Find.Execute("75+11=86", $true, $false, $false, $false, $false, $true, 1, $false, "37+34=71", 2) | Out-Null
$d.Content.Find.Execute("45+2=47", $true, $false, $false, $false, $false, $true, 1, $false, "52-47=5", 2) | Out-Null
$d.Content.Find.Execute("74-1=73", $true, $false, $false, $false, $false, $true, 1, $false, "86-61=25", 2) | Out-Null
$d.Content.Find.Execute("20+42=62", $true, $false, $false, $false, $false, $true, 1, $false, "26+18=44", 2) | Out-Null
$d.Content.Find.Execute("71-34=37", $true, $false, $false, $false, $false, $true, 1, $false, "4+21=25", 2) | Out-Null
$d.Content.Find.Execute("46-33=13", $true, $false, $false, $false, $false, $true, 1, $false, "20+77=97", 2) | Out-Null
$d.Content.Find.Execute("36-19=17", $true, $false, $false, $false, $false, $true, 1, $false, "61-42=19", 2) | Out-Null
$d.Content.Find.Execute("2+11=13", $true, $false, $false, $false, $false, $true, 1, $false, "48-41=7", 2) | Out-Null
$d.Content.Find.Execute("65+29=94", $true, $false, $false, $false, $false, $true, 1, $false, "83-24=59", 2) | Out-Null
$d.Content.Find.Execute("9+29=38", $true, $false, $false, $false, $false, $true, 1, $false, "69-37=32", 2) | Out-Null
$d.Content.Find.Execute("66+2=68", $true, $false, $false, $false, $false, $true, 1, $false, "2+81=83", 2) | Out-Null
$d.Content.Find.Execute("20+46=66", $true, $false, $false, $false, $false, $true, 1, $false, "3+88=91", 2) | Out-Null
$d.Content.Find.Execute("8+10=18", $true, $false, $false, $false, $false, $true, 1, $false, "60-8=52", 2) | Out-Null
$d.Content.Find.Execute("95-6=89", $true, $false, $false, $false, $false, $true, 1, $false, "74-19=55", 2) | Out-Null
$d.Content.Find.Execute("42+38=80", $true, $false, $false, $false, $false, $true, 1, $false, "43+50=93", 2) | Out-Null
$d.Content.Find.Execute("21+65=86", $true, $false, $false, $false, $false, $true, 1, $false, "58+35=93", 2) | Out-Null
$d.Content.Find.Execute("92-28=64", $true, $false, $false, $false, $false, $true, 1, $false, "14+9=23", 2) | Out-Null
$d.Content.Find.Execute("62-26=36", $true, $false, $false, $false, $false, $true, 1, $false, "44+32=76", 2) | Out-Null
$d.Content.Find.Execute("5+65=70", $true, $false, $false, $false, $false, $true, 1, $false, "47-34=13", 2) | Out-Null
$d.Content.Find.Execute("81-43=38", $true, $false, $false, $false, $false, $true, 1, $false, "27+61=88", 2) | Out-Null
$d.Content.Find.Execute("1+58=59", $true, $false, $false, $false, $false, $true, 1, $false, "17+20=37", 2) | Out-Null
$d.Content.Find.Execute("88-68=20", $true, $false, $false, $false, $false, $true, 1, $false, "4+75=79", 2) | Out-Null
$d.Content.Find.Execute("69-68=1", $true, $false, $false, $false, $false, $true, 1, $false, "31-14=17", 2) | Out-Null
$d.Content.Find.Execute("30+18=48", $true, $false, $false, $false, $false, $true, 1, $false, "83+5=88", 2) | Out-Null
$d.Content.Find.Execute("46-2=44", $true, $false, $false, $false, $false, $true, 1, $false, "94-93=1", 2) | Out-Null
$d.Content.Find.Execute("13+78=91", $true, $false, $false, $false, $false, $true, 1, $false, "84-54=30", 2) | Out-Null
$d.Content.Find.Execute("99-25=74", $true, $false, $false, $false, $false, $true, 1, $false, "2+46=48", 2) | Out-Null
$d.Content.Find.Execute("6+86=92", $true, $false, $false, $false, $false, $true, 1, $false, "83-24=59", 2) | Out-Null
$d.Content.Find.Execute("78+8=86", $true, $false, $false, $false, $false, $true, 1, $false, "13+24=37", 2) | Out-Null
$d.Content.Find.Execute("54-3=51", $true, $false, $false, $false, $false, $true, 1, $false, "6+70=76", 2) | Out-Null
$d.Content.Find.Execute("57-28=29", $true, $false, $false, $false, $false, $true, 1, $false, "49-34=15", 2) | Out-Null
$d.Content.Find.Execute("60-13=47", $true, $false, $false, $false, $false, $true, 1, $false, "34-1=33", 2) | Out-Null
$d.Content.Find.Execute("87-18=69", $true, $false, $false, $false, $false, $true, 1, $false, "76-14=62", 2) | Out-Null
$d.Content.Find.Execute("72-54=18", $true, $false, $false, $false, $false, $true, 1, $false, "79+5=84", 2) | Out-Null
$d.Content.Find.Execute("20+64=84", $true, $false, $false, $false, $false, $true, 1, $false, "54+19=73", 2) | Out-Null
$d.Content.Find.Execute("23+33=56", $true, $false, $false, $false, $false, $true, 1, $false, "68-38=30", 2) | Out-Null
$d.Content.Find.Execute("97-4=93", $true, $false, $false, $false, $false, $true, 1, $false, "69-46=23", 2) | Out-Null
$d.Content.Find.Execute("58-42=16", $true, $false, $false, $false, $false, $true, 1, $false, "17+4=21", 2) | Out-Null
$d.Content.Find.Execute("11+28=39", $true, $false, $false, $false, $false, $true, 1, $false, "35+21=56", 2) | Out-Null
$d.Content.Find.Execute("96-25=71", $true, $false, $false, $false, $false, $true, 1, $false, "13+56=69", 2) | Out-Null
$d.Content.Find.Execute("2+57=59", $true, $false, $false, $false, $false, $true, 1, $false, "33+35=68", 2) | Out-Null
$d.Content.Find.Execute("91-0=91", $true, $false, $false, $false, $false, $true, 1, $false, "43+43=86", 2) | Out-Null
$d.Content.Find.Execute("52+33=85", $true, $false, $false, $false, $false, $true, 1, $false, "16+11=27", 2) | Out-Null
$d.Content.Find.Execute("14+29=43", $true, $false, $false, $false, $false, $true, 1, $false, "71+16=87", 2) | Out-Null
$d.Content.Find.Execute("44-25=19", $true, $false, $false, $false, $false, $true, 1, $false, "5+55=60", 2) | Out-Null
$d.Content.Find.Execute("28+13=41", $true, $false, $false, $false, $false, $true, 1, $false, "1+39=40", 2) | Out-Null
$d.Content.Find.Execute("46+45=91", $true, $false, $false, $false, $false, $true, 1, $false, "24-12=12", 2) | Out-Null
$d.Content.Find.Execute("93-13=80", $true, $false, $false, $false, $false, $true, 1, $false, "51+25=76", 2) | Out-Null
$d.Content.Find.Execute("4+70=74", $true, $false, $false, $false, $false, $true, 1, $false, "32+20=52", 2) | Out-Null
$d.Content.Find.Execute("78-48=30", $true, $false, $false, $false, $false, $true, 1, $false, "27+2=29", 2) | Out-Null
$d.Content.Find.Execute("32+56=88", $true, $false, $false, $false, $false, $true, 1, $false, "61-29=32", 2) | Out-Null
$d.Content.Find.Execute("14+0=14", $true, $false, $false, $false, $false, $true, 1, $false, "96-75=21", 2) | Out-Null
$d.Content.Find.Execute("23+72=95", $true, $false, $false, $false, $false, $true, 1, $false, "10+16=26", 2) | Out-Null
$d.Content.Find.Execute("71-32=39", $true, $false, $false, $false, $false, $true, 1, $false, "41-5=36", 2) | Out-Null
$d.Content.Find.Execute("25-10=15", $true, $false, $false, $false, $false, $true, 1, $false, "90-8=82", 2) | Out-Null
$d.Content.Find.Execute("4+8=12", $true, $false, $false, $false, $false, $true, 1, $false, "52+38=90", 2) | Out-Null
$d.Content.Find.Execute("89-11=78", $true, $false, $false, $false, $false, $true, 1, $false, "3+39=42", 2) | Out-Null
$d.Content.Find.Execute("65+25=90", $true, $false, $false, $false, $false, $true, 1, $false, "84-74=10", 2) | Out-Null
$d.Content.Find.Execute("56+7=63", $true, $false, $false, $false, $false, $true, 1, $false, "24+56=80", 2) | Out-Null
$d.Content.Find.Execute("28-5=23", $true, $false, $false, $false, $false, $true, 1, $false, "76+7=83", 2) | Out-Null
$d.Content.Find.Execute("20+61=81", $true, $false, $false, $false, $false, $true, 1, $false, "83+16=99", 2) | Out-Null
$d.Content.Find.Execute("43+25=68", $true, $false, $false, $false, $false, $true, 1, $false, "27-14=13", 2) | Out-Null
$d.Content.Find.Execute("58-46=12", $true, $false, $false, $false, $false, $true, 1, $false, "7-7=0", 2) | Out-Null
$d.Content.Find.Execute("61+8=69", $true, $false, $false, $false, $false, $true, 1, $false, "48-9=39", 2) | Out-Null
$d.Content.Find.Execute("6+37=43", $true, $false, $false, $false, $false, $true, 1, $false, "3+47=50", 2) | Out-Null
$d.Content.Find.Execute("63+29=92", $true, $false, $false, $false, $false, $true, 1, $false, "73-62=11", 2) | Out-Null
$d.Content.Find.Execute("99-1=98", $true, $false, $false, $false, $false, $true, 1, $false, "43+51=94", 2) | Out-Null
$d.Content.Find.Execute("38-6=32", $true, $false, $false, $false, $false, $true, 1, $false, "11+25=36", 2) | Out-Null
$d.Content.Find.Execute("11+81=92", $true, $false, $false, $false, $false, $true, 1, $false, "22-15=7", 2) | Out-Null
$d.Content.Find.Execute("90-89=1", $true, $false, $false, $false, $false, $true, 1, $false, "28+3=31", 2) | Out-Null
$d.Content.Find.Execute("30+41=71", $true, $false, $false, $false, $false, $true, 1, $false, "97-46=51", 2) | Out-Null
$d.Content.Find.Execute("77-44=33", $true, $false, $false, $false, $false, $true, 1, $false, "14+11=25", 2) | Out-Null
$d.Content.Find.Execute("60-47=13", $true, $false, $false, $false, $false, $true, 1, $false, "90-64=26", 2) | Out-Null
$d.Content.Find.Execute("64+26=90", $true, $false, $false, $false, $false, $true, 1, $false, "59-5=54", 2) | Out-Null
$d.Content.Find.Execute("2+20=22", $true, $false, $false, $false, $false, $true, 1, $false, "29+37=66", 2) | Out-Null
$d.Content.Find.Execute("80-22=58", $true, $false, $false, $false, $false, $true, 1, $false, "63-27=36", 2) | Out-Null
$d.Content.Find.Execute("2+41=43", $true, $false, $false, $false, $false, $true, 1, $false, "82+8=90", 2) | Out-Null
$d.Content.Find.Execute("16+43=59", $true, $false, $false, $false, $false, $true, 1, $false, "0+99=99", 2) | Out-Null
$d.Content.Find.Execute("94-55=39", $true, $false, $false, $false, $false, $true, 1, $false, "24-20=4", 2) | Out-Null
$d.Content.Find.Execute("75+4=79", $true, $false, $false, $false, $false, $true, 1, $false, "28-16=12", 2) | Out-Null
$d.Content.Find.Execute("19+8=27", $true, $false, $false, $false, $false, $true, 1, $false, "52-19=33", 2) | Out-Null
$d.Content.Find.Execute("57+3=60", $true, $false, $false, $false, $false, $true, 1, $false, "15+26=41", 2) | Out-Null
$d.Content.Find.Execute("95-74=21", $true, $false, $false, $false, $false, $true, 1, $false, "68-59=9", 2) | Out-Null
$d.Content.Find.Execute("57+42=99", $true, $false, $false, $false, $false, $true, 1, $false, "35-5=30", 2) | Out-Null
$d.Content.Find.Execute("34+5=39", $true, $false, $false, $false, $false, $true, 1, $false, "6+41=47", 2) | Out-Null
$d.Content.Find.Execute("26+25=51", $true, $false, $false, $false, $false, $true, 1, $false, "59-3=56", 2) | Out-Null
$d.Content.Find.Execute("5+15=20", $true, $false, $false, $false, $false, $true, 1, $false, "33-0=33", 2) | Out-Null
$d.Content.Find.Execute("96-32=64", $true, $false, $false, $false, $false, $true, 1, $false, "11+45=56", 2) | Out-Null
$d.Content.Find.Execute("56+13=69", $true, $false, $false, $false, $false, $true, 1, $false, "12+84=96", 2) | Out-Null
$d.Content.Find.Execute("58-32=26", $true, $false, $false, $false, $false, $true, 1, $false, "96-8=88", 2) | Out-Null
$d.Content.Find.Execute("22+67=89", $true, $false, $false, $false, $false, $true, 1, $false, "83-0=83", 2) | Out-Null
$d.Content.Find.Execute("73-25=48", $true, $false, $false, $false, $false, $true, 1, $false, "57+25=82", 2) | Out-Null
$d.Content.Find.Execute("18+26=44", $true, $false, $false, $false, $false, $true, 1, $false, "68-60=8", 2) | Out-Null
$d.Content.Find.Execute("55+40=95", $true, $false, $false, $false, $false, $true, 1, $false, "39+9=48", 2) | Out-Null
$d.Content.Find.Execute("91-32=59", $true, $false, $false, $false, $false, $true, 1, $false, "1+18=19", 2) | Out-Null
